# This script reproduces the commit "Data importation and start of main code":
# it consolidates the three existing monthly/yearly grids (Numero spettacoli,
# Ingressi, Spesa del pubblico) into one long/tidy table on a new "Sheet1" tab
# at the end of the workbook, with a Date column (month-start, mmm-yy format).

$wb = $excel.ActiveWorkbook

# --- Add the new worksheet as the last tab in the workbook ---
$sheetCount = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($sheetCount)
$new = $wb.Worksheets.Add($null, $lastSheet)

# --- Header row ---
$new.Cells.Item(1, 1).Value = "Date"
$new.Cells.Item(1, 2).Value = "Numero spettacoli"
$new.Cells.Item(1, 3).Value = "Ingressi"
$new.Cells.Item(1, 4).Value = "Spesa del pubblico"

# --- Data rows: one row per (year, month) combination, Jan 2018 .. Dec 2022 ---
$new.Cells.Item(2, 1).Value = 43101
$new.Cells.Item(2, 2).Value = 319556
$new.Cells.Item(2, 3).Value = 12815316
$new.Cells.Item(2, 4).Value = 99527809.3999999
$new.Cells.Item(3, 1).Value = 43132
$new.Cells.Item(3, 2).Value = 279211
$new.Cells.Item(3, 3).Value = 10586422
$new.Cells.Item(3, 4).Value = 78156505.44999997
$new.Cells.Item(4, 1).Value = 43160
$new.Cells.Item(4, 2).Value = 302038
$new.Cells.Item(4, 3).Value = 8408634
$new.Cells.Item(4, 4).Value = 60683989.14
$new.Cells.Item(5, 1).Value = 43191
$new.Cells.Item(5, 2).Value = 287672
$new.Cells.Item(5, 3).Value = 7429757
$new.Cells.Item(5, 4).Value = 52284212.71000002
$new.Cells.Item(6, 1).Value = 43221
$new.Cells.Item(6, 2).Value = 262643
$new.Cells.Item(6, 3).Value = 5579501
$new.Cells.Item(6, 4).Value = 45098825.43000006
$new.Cells.Item(7, 1).Value = 43252
$new.Cells.Item(7, 2).Value = 212699
$new.Cells.Item(7, 3).Value = 3675420
$new.Cells.Item(7, 4).Value = 28124127.30000001
$new.Cells.Item(8, 1).Value = 43282
$new.Cells.Item(8, 2).Value = 176592
$new.Cells.Item(8, 3).Value = 2972751
$new.Cells.Item(8, 4).Value = 19073912.039999984
$new.Cells.Item(9, 1).Value = 43313
$new.Cells.Item(9, 2).Value = 187101
$new.Cells.Item(9, 3).Value = 5110096
$new.Cells.Item(9, 4).Value = 37118997.81000001
$new.Cells.Item(10, 1).Value = 43344
$new.Cells.Item(10, 2).Value = 264158
$new.Cells.Item(10, 3).Value = 6789039
$new.Cells.Item(10, 4).Value = 52631424.29999997
$new.Cells.Item(11, 1).Value = 43374
$new.Cells.Item(11, 2).Value = 281824
$new.Cells.Item(11, 3).Value = 7709847
$new.Cells.Item(11, 4).Value = 57354454.720000006
$new.Cells.Item(12, 1).Value = 43405
$new.Cells.Item(12, 2).Value = 279071
$new.Cells.Item(12, 3).Value = 8672730
$new.Cells.Item(12, 4).Value = 62379808.1599999
$new.Cells.Item(13, 1).Value = 43435
$new.Cells.Item(13, 2).Value = 312125
$new.Cells.Item(13, 3).Value = 11938186
$new.Cells.Item(13, 4).Value = 89976686.61999996
$new.Cells.Item(14, 1).Value = 43466
$new.Cells.Item(14, 2).Value = 308748
$new.Cells.Item(14, 3).Value = 12489100
$new.Cells.Item(14, 4).Value = 100171903.29999997
$new.Cells.Item(15, 1).Value = 43497
$new.Cells.Item(15, 2).Value = 272951
$new.Cells.Item(15, 3).Value = 8467518
$new.Cells.Item(15, 4).Value = 61659419.87999999
$new.Cells.Item(16, 1).Value = 43525
$new.Cells.Item(16, 2).Value = 301053
$new.Cells.Item(16, 3).Value = 8068654
$new.Cells.Item(16, 4).Value = 59014191.74
$new.Cells.Item(17, 1).Value = 43556
$new.Cells.Item(17, 2).Value = 288046
$new.Cells.Item(17, 3).Value = 10290518
$new.Cells.Item(17, 4).Value = 77990289.28
$new.Cells.Item(18, 1).Value = 43586
$new.Cells.Item(18, 2).Value = 268636
$new.Cells.Item(18, 3).Value = 7796085
$new.Cells.Item(18, 4).Value = 59729552.51000001
$new.Cells.Item(19, 1).Value = 43617
$new.Cells.Item(19, 2).Value = 237477
$new.Cells.Item(19, 3).Value = 4286370
$new.Cells.Item(19, 4).Value = 31849329.549999993
$new.Cells.Item(20, 1).Value = 43647
$new.Cells.Item(20, 2).Value = 207763
$new.Cells.Item(20, 3).Value = 4909151
$new.Cells.Item(20, 4).Value = 37861308.36
$new.Cells.Item(21, 1).Value = 43678
$new.Cells.Item(21, 2).Value = 210585
$new.Cells.Item(21, 3).Value = 7051071
$new.Cells.Item(21, 4).Value = 54350700.96
$new.Cells.Item(22, 1).Value = 43709
$new.Cells.Item(22, 2).Value = 252819
$new.Cells.Item(22, 3).Value = 7819079
$new.Cells.Item(22, 4).Value = 62794563.51
$new.Cells.Item(23, 1).Value = 43739
$new.Cells.Item(23, 2).Value = 284646
$new.Cells.Item(23, 3).Value = 9632400
$new.Cells.Item(23, 4).Value = 72910605.82000001
$new.Cells.Item(24, 1).Value = 43770
$new.Cells.Item(24, 2).Value = 291507
$new.Cells.Item(24, 3).Value = 9937379
$new.Cells.Item(24, 4).Value = 70529036.69000003
$new.Cells.Item(25, 1).Value = 43800
$new.Cells.Item(25, 2).Value = 321114
$new.Cells.Item(25, 3).Value = 13691888
$new.Cells.Item(25, 4).Value = 103129044.46
$new.Cells.Item(26, 1).Value = 43831
$new.Cells.Item(26, 2).Value = 342511
$new.Cells.Item(26, 3).Value = 16606889
$new.Cells.Item(26, 4).Value = 121440619.01000008
$new.Cells.Item(27, 1).Value = 43862
$new.Cells.Item(27, 2).Value = 251640
$new.Cells.Item(27, 3).Value = 7080523
$new.Cells.Item(27, 4).Value = 48923544.429999985
$new.Cells.Item(28, 1).Value = 43891
$new.Cells.Item(28, 2).Value = 14247
$new.Cells.Item(28, 3).Value = 217518
$new.Cells.Item(28, 4).Value = 1440488.9799999997
$new.Cells.Item(29, 1).Value = 43922
$new.Cells.Item(29, 4).Value = 325.5
$new.Cells.Item(30, 1).Value = 43952
$new.Cells.Item(30, 4).Value = 198
$new.Cells.Item(31, 1).Value = 43983
$new.Cells.Item(31, 2).Value = 6944
$new.Cells.Item(31, 3).Value = 111677
$new.Cells.Item(31, 4).Value = 685211.22
$new.Cells.Item(32, 1).Value = 44013
$new.Cells.Item(32, 2).Value = 23934
$new.Cells.Item(32, 3).Value = 658929
$new.Cells.Item(32, 4).Value = 3309540.71
$new.Cells.Item(33, 1).Value = 44044
$new.Cells.Item(33, 2).Value = 59427
$new.Cells.Item(33, 3).Value = 1381564
$new.Cells.Item(33, 4).Value = 7665196.540000002
$new.Cells.Item(34, 1).Value = 44075
$new.Cells.Item(34, 2).Value = 158221
$new.Cells.Item(34, 3).Value = 2443035
$new.Cells.Item(34, 4).Value = 17362260.849999994
$new.Cells.Item(35, 1).Value = 44105
$new.Cells.Item(35, 2).Value = 139085
$new.Cells.Item(35, 3).Value = 1805636
$new.Cells.Item(35, 4).Value = 11266608.920000006
$new.Cells.Item(36, 1).Value = 44136
$new.Cells.Item(36, 4).Value = 4823.5
$new.Cells.Item(37, 1).Value = 44166
$new.Cells.Item(37, 4).Value = 56332
$new.Cells.Item(38, 1).Value = 44197
$new.Cells.Item(38, 4).Value = 27483
$new.Cells.Item(39, 1).Value = 44228
$new.Cells.Item(39, 4).Value = 14353.5
$new.Cells.Item(40, 1).Value = 44256
$new.Cells.Item(40, 4).Value = 2164.8
$new.Cells.Item(41, 1).Value = 44287
$new.Cells.Item(41, 2).Value = 2174
$new.Cells.Item(41, 3).Value = 44478
$new.Cells.Item(41, 4).Value = 289498.68
$new.Cells.Item(42, 1).Value = 44317
$new.Cells.Item(42, 2).Value = 61019
$new.Cells.Item(42, 3).Value = 801049
$new.Cells.Item(42, 4).Value = 5164938.140000001
$new.Cells.Item(43, 1).Value = 44348
$new.Cells.Item(43, 2).Value = 124073
$new.Cells.Item(43, 3).Value = 1415528
$new.Cells.Item(43, 4).Value = 9102229.35
$new.Cells.Item(44, 1).Value = 44378
$new.Cells.Item(44, 2).Value = 137767
$new.Cells.Item(44, 3).Value = 2318804
$new.Cells.Item(44, 4).Value = 15407893.210000005
$new.Cells.Item(45, 1).Value = 44409
$new.Cells.Item(45, 2).Value = 156915
$new.Cells.Item(45, 3).Value = 2679886
$new.Cells.Item(45, 4).Value = 18074571.25
$new.Cells.Item(46, 1).Value = 44440
$new.Cells.Item(46, 2).Value = 187832
$new.Cells.Item(46, 3).Value = 3272931
$new.Cells.Item(46, 4).Value = 22875671.26999999
$new.Cells.Item(47, 1).Value = 44470
$new.Cells.Item(47, 2).Value = 217862
$new.Cells.Item(47, 3).Value = 5048158
$new.Cells.Item(47, 4).Value = 35585479.21
$new.Cells.Item(48, 1).Value = 44501
$new.Cells.Item(48, 2).Value = 225776
$new.Cells.Item(48, 3).Value = 5050651
$new.Cells.Item(48, 4).Value = 34525123.07999999
$new.Cells.Item(49, 1).Value = 44531
$new.Cells.Item(49, 2).Value = 237348
$new.Cells.Item(49, 3).Value = 6101141
$new.Cells.Item(49, 4).Value = 44847154.66
$new.Cells.Item(50, 1).Value = 44562
$new.Cells.Item(50, 2).Value = 250091
$new.Cells.Item(50, 3).Value = 4536594
$new.Cells.Item(50, 4).Value = 32381225.9
$new.Cells.Item(51, 1).Value = 44593
$new.Cells.Item(51, 2).Value = 170035
$new.Cells.Item(51, 3).Value = 2978630
$new.Cells.Item(51, 4).Value = 20299393.849999998
$new.Cells.Item(52, 1).Value = 44621
$new.Cells.Item(52, 2).Value = 199764
$new.Cells.Item(52, 3).Value = 3658086
$new.Cells.Item(52, 4).Value = 24659053.099999998
$new.Cells.Item(53, 1).Value = 44652
$new.Cells.Item(53, 2).Value = 212027
$new.Cells.Item(53, 3).Value = 4048463
$new.Cells.Item(53, 4).Value = 27809817.410000008
$new.Cells.Item(54, 1).Value = 44682
$new.Cells.Item(54, 2).Value = 187349
$new.Cells.Item(54, 3).Value = 3792309
$new.Cells.Item(54, 4).Value = 27234239.39000001
$new.Cells.Item(55, 1).Value = 44713
$new.Cells.Item(55, 2).Value = 159464
$new.Cells.Item(55, 3).Value = 3065783
$new.Cells.Item(55, 4).Value = 21926129.69000001
$new.Cells.Item(56, 1).Value = 44743
$new.Cells.Item(56, 2).Value = 125896
$new.Cells.Item(56, 3).Value = 2687118
$new.Cells.Item(56, 4).Value = 18478583.930000003
$new.Cells.Item(57, 1).Value = 44774
$new.Cells.Item(57, 2).Value = 126277
$new.Cells.Item(57, 3).Value = 3036152
$new.Cells.Item(57, 4).Value = 20729046.200000003
$new.Cells.Item(58, 1).Value = 44805
$new.Cells.Item(58, 2).Value = 178988
$new.Cells.Item(58, 3).Value = 3667427
$new.Cells.Item(58, 4).Value = 22258833.90000001
$new.Cells.Item(59, 1).Value = 44835
$new.Cells.Item(59, 2).Value = 211458
$new.Cells.Item(59, 3).Value = 4364739
$new.Cells.Item(59, 4).Value = 30259469.300000016
$new.Cells.Item(60, 1).Value = 44866
$new.Cells.Item(60, 2).Value = 208221
$new.Cells.Item(60, 3).Value = 4946407
$new.Cells.Item(60, 4).Value = 33332163.960000012
$new.Cells.Item(61, 1).Value = 44896
$new.Cells.Item(61, 2).Value = 225459
$new.Cells.Item(61, 3).Value = 6869912
$new.Cells.Item(61, 4).Value = 53667485.470000006

# --- Format the Date column as "mmm-yy" (numFmtId 17) ---
$new.Range("A2:A61").NumberFormat = "mmm-yy"

# --- Restore/refresh selections on the pre-existing sheets ---
$ws1 = $wb.Worksheets.Item("Numero spettacoli")
$ws1.Activate() | Out-Null
$ws1.Range("B2:F13").Select() | Out-Null

$ws2 = $wb.Worksheets.Item("Ingressi")
$ws2.Activate() | Out-Null
$ws2.Range("B2:F13").Select() | Out-Null

$ws3 = $wb.Worksheets.Item("Spesa del pubblico")
$ws3.Activate() | Out-Null
$ws3.Range("B2:F13").Select() | Out-Null

# --- Finish with the new sheet active/selected (matches the saved workbook state) ---
$new.Activate() | Out-Null
$new.Range("G33").Select() | Out-Null
